$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.495.60"
$ws.Range("E2").Value = "  +0.20%  "
$ws.Range("D3").Value = "1.836.89"
$ws.Range("E3").Value = "  -0.29%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "258.08"
$ws.Range("E5").Value = "  -1.16%  "
$ws.Range("E6").Value = "  +0.04%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5238"
$ws.Range("E7").Value = "  +0.50%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3149"
$ws.Range("E8").Value = "  -3.92%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06782"
$ws.Range("E9").Value = "  -0.04%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.57"
$ws.Range("E10").Value = "  -0.32%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.7741"
$ws.Range("E11").Value = "  +0.14%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07748"
$ws.Range("D13").Value = "1.830.33"
$ws.Range("E13").Value = "  -0.72%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "87.50"
$ws.Range("E14").Value = "  -0.80%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.994"
$ws.Range("E15").Value = "  -0.65%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.000"
$ws.Range("E16").Value = "  +0.02%  "
$ws.Range("E18").Value = "  +0.07%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007916"
$ws.Range("E19").Value = "  -0.49%  "
$ws.Range("D20").Value = "26.525.04"
$ws.Range("E20").Value = "  +0.17%  "
$ws.Range("D21").Value = "2.066.99"
$ws.Range("E21").Value = "  -0.21%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.579"
$ws.Range("E22").Value = "  -0.02%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.942"
$ws.Range("E23").Value = "  -0.57%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.274"
$ws.Range("E24").Value = "  -2.25%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "142.17"
$ws.Range("E25").Value = "  -1.59%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.195"
$ws.Range("E26").Value = "  -0.48%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.671"
$ws.Range("E27").Value = "  +0.93%  "
$ws.Range("E28").Value = "  -0.76%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "110.49"
$ws.Range("E29").Value = "  -0.80%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.145"
$ws.Range("E30").Value = "  -1.18%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08709"
$ws.Range("E31").Value = "  -0.12%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.045"
$ws.Range("E32").Value = "  -2.17%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04839"
$ws.Range("E33").Value = "  +0.86%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.133"
$ws.Range("E34").Value = "  +0.28%  "
$ws.Range("B35").Value = "HuobiToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.861"
$ws.Range("E35").Value = "  +0.87%  "
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7150"
$ws.Range("E36").Value = "  +0.65%  "
$ws.Range("E37").Value = "  -0.47%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.213"
$ws.Range("E38").Value = "  -1.07%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01721"
$ws.Range("E39").Value = "  -2.34%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.4787"
$ws.Range("E40").Value = "  -1.16%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8908"
$ws.Range("E41").Value = "  -0.27%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "109.40"
$ws.Range("E42").Value = "  -2.31%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.902"
$ws.Range("E43").Value = "  -2.76%  "
$ws.Range("E44").Value = "  +0.09%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.601"
$ws.Range("E45").Value = "  -1.71%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4134"
$ws.Range("E46").Value = "  -0.87%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.928"
$ws.Range("E47").Value = "  -0.28%  "
$ws.Range("E48").Value = "  -1.15%  "
$ws.Range("E49").Value = "  +0.53%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "34.58"
$ws.Range("E50").Value = "  -1.09%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.8907"
$ws.Range("E51").Value = "  +0.11%  "
